{"js": "// Locate the last BodyText paragraph (\"This one comes from the deploy\n// keys\"), rewrite its text, then append two more BodyText paragraphs\n// right after it (still inside the \"hello\" bookmark range).\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst oldText = \"This one comes from the deploy keys\";\nconst target = paragraphs.items.find((p) => p.text === oldText);\nif (!target) {\n  throw new Error(\"Could not find paragraph: \" + oldText);\n}\n\n// 1. Replace the existing run text in place.\ntarget.insertText(\"Another test just in case\", \"Replace\");\n\n// 2. Insert the two new paragraphs right after it, matching the\n//    surrounding \"Body Text\" style.\nconst p1 = target.insertParagraph(\n  \"Another test with the right branch.\",\n  \"After\"\n);\np1.style = \"Body Text\";\n\nconst p2 = p1.insertParagraph(\"Fix / else / then is tested\", \"After\");\np2.style = \"Body Text\";\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1. Replace the text of the last \"BodyText\" paragraph in place.\n$find = $d.Content.Find\n$find.Text = \"This one comes from the deploy keys\"\n$find.Replacement.Text = \"Another test just in case\"\n$find.Execute(\n    $find.Text,\n    $true, $false, $false, $false, $false,\n    $true, 0, $false,\n    $find.Replacement.Text,\n    2\n)\n\n# 2. Append two more \"Body Text\" paragraphs right after it (still before\n#    the bookmark end), matching the style of their neighbours.\n$target = $d.Paragraphs.Last\n$target.Range.InsertParagraphAfter()\n\n$p1 = $d.Paragraphs.Last\n$p1.Style = $target.Style\n$p1.Range.Text = \"Another test with the right branch.\"\n\n$p1.Range.InsertParagraphAfter()\n\n$p2 = $d.Paragraphs.Last\n$p2.Style = $target.Style\n$p2.Range.Text = \"Fix / else / then is tested\"\n"}
